# Generate Report for Handoff
# Updates the localization-status report after the source markdown file was
# regenerated under a new GUID-named file (84b9a225-... -> 9c9a7992-...),
# refreshed the handoff timestamps/xliff names, and reset the (not yet
# produced) handback columns for both the zh-cn and de-de target sheets.

$wb = $excel.ActiveWorkbook

$oldGuidFile = "84b9a225-6938-4a42-9132-98900a48dad2.md"
$newGuidFile = "9c9a7992-a529-4b1d-b21f-3df28d1196f7.md"

# ----------------------------------------------------------------------
# Overview sheet
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuidFile
$wsOverview.Range("B2").Value = "e2e\" + $newGuidFile
$wsOverview.Range("G2").Value = "2016-08-18 19:03:48"

# ----------------------------------------------------------------------
# zh-cn sheet
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newGuidFile
$wsZhCn.Range("G2").Value = "9c9a7992-a529-4b1d-b21f-3df28d1196f7.ab10d89ea9dc9b37266dd1954d819e1017e8bac8.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-18 19:03:42"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# ----------------------------------------------------------------------
# de-de sheet
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newGuidFile
$wsDeDe.Range("G2").Value = "9c9a7992-a529-4b1d-b21f-3df28d1196f7.ab10d89ea9dc9b37266dd1954d819e1017e8bac8.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-18 19:03:48"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"
